# BIS-1002: removed "Internal Assignment" column from export.
# The "Internal Assignment" header (O4) and its "FALSE" values (O5:O7)
# are cleared out, while keeping the existing cell formatting/styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O4:O7").ClearContents()
